# Update cryptos list data (prices and volume %) scraped on
# Sun Mar 26 14:00:15 UTC 2023, and fix the swapped HuobiToken /
# Filecoin rows (33 and 34) so coin/link/price/volume line up again.
#
# Values are prefixed with a leading apostrophe so Excel stores them
# as literal text (matching the source workbook's inline-string
# cells) instead of auto-converting numeric-looking strings like
# "0.9994" into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''28.105.52'
$ws.Cells.Item(2, 5).Value = '''  +1.93%  '
$ws.Cells.Item(3, 4).Value = '''1.792.16'
$ws.Cells.Item(3, 5).Value = '''  +2.14%  '
$ws.Cells.Item(4, 4).Value = '''0.9994'
$ws.Cells.Item(4, 5).Value = '''  -0.08%  '
$ws.Cells.Item(5, 4).Value = '''327.26'
$ws.Cells.Item(5, 5).Value = '''  +0.89%  '
$ws.Cells.Item(6, 4).Value = '''0.9982'
$ws.Cells.Item(7, 4).Value = '''0.4527'
$ws.Cells.Item(7, 5).Value = '''  +1.33%  '
$ws.Cells.Item(8, 5).Value = '''  +0.53%  '
$ws.Cells.Item(9, 4).Value = '''0.07522'
$ws.Cells.Item(9, 5).Value = '''  +0.39%  '
$ws.Cells.Item(10, 4).Value = '''42.54'
$ws.Cells.Item(10, 5).Value = '''  +1.15%  '
$ws.Cells.Item(11, 4).Value = '''1.115'
$ws.Cells.Item(11, 5).Value = '''  +1.90%  '
$ws.Cells.Item(12, 4).Value = '''0.9986'
$ws.Cells.Item(12, 5).Value = '''  -0.14%  '
$ws.Cells.Item(13, 4).Value = '''21.09'
$ws.Cells.Item(13, 5).Value = '''  +1.18%  '
$ws.Cells.Item(14, 4).Value = '''6.088'
$ws.Cells.Item(14, 5).Value = '''  +0.96%  '
$ws.Cells.Item(15, 4).Value = '''7.270'
$ws.Cells.Item(15, 5).Value = '''  +2.12%  '
$ws.Cells.Item(16, 4).Value = '''1.784.43'
$ws.Cells.Item(16, 5).Value = '''  +2.44%  '
$ws.Cells.Item(17, 4).Value = '''94.34'
$ws.Cells.Item(17, 5).Value = '''  +1.17%  '
$ws.Cells.Item(18, 4).Value = '''0.00001066'
$ws.Cells.Item(18, 5).Value = '''  +0.39%  '
$ws.Cells.Item(19, 4).Value = '''0.06456'
$ws.Cells.Item(19, 5).Value = '''  +0.64%  '
$ws.Cells.Item(20, 4).Value = '''0.9983'
$ws.Cells.Item(20, 5).Value = '''  -0.09%  '
$ws.Cells.Item(21, 4).Value = '''17.25'
$ws.Cells.Item(21, 5).Value = '''  +2.38%  '
$ws.Cells.Item(22, 4).Value = '''5.844'
$ws.Cells.Item(22, 5).Value = '''  +0.50%  '
$ws.Cells.Item(23, 4).Value = '''28.120.41'
$ws.Cells.Item(23, 5).Value = '''  +1.80%  '
$ws.Cells.Item(24, 4).Value = '''11.42'
$ws.Cells.Item(25, 4).Value = '''2.089'
$ws.Cells.Item(25, 5).Value = '''  -1.01%  '
$ws.Cells.Item(26, 4).Value = '''163.75'
$ws.Cells.Item(26, 5).Value = '''  +0.67%  '
$ws.Cells.Item(27, 4).Value = '''20.42'
$ws.Cells.Item(27, 5).Value = '''  -0.08%  '
$ws.Cells.Item(28, 4).Value = '''1.989.42'
$ws.Cells.Item(28, 5).Value = '''  +2.17%  '
$ws.Cells.Item(29, 4).Value = '''2.278'
$ws.Cells.Item(29, 5).Value = '''  +9.35%  '
$ws.Cells.Item(30, 4).Value = '''126.69'
$ws.Cells.Item(30, 5).Value = '''  +0.25%  '
$ws.Cells.Item(31, 4).Value = '''1.122'
$ws.Cells.Item(31, 5).Value = '''  +3.82%  '
$ws.Cells.Item(32, 4).Value = '''0.09202'
$ws.Cells.Item(32, 5).Value = '''  +1.28%  '
$ws.Cells.Item(33, 2).Value = '''Filecoin'
$ws.Cells.Item(33, 3).Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).Value = '''5.621'
$ws.Cells.Item(33, 5).Value = '''  +1.52%  '
$ws.Cells.Item(34, 2).Value = '''HuobiToken'
$ws.Cells.Item(34, 3).Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(34, 4).Value = '''3.678'
$ws.Cells.Item(34, 5).Value = '''  +0.26%  '
$ws.Cells.Item(35, 4).Value = '''11.98'
$ws.Cells.Item(35, 5).Value = '''  +0.07%  '
$ws.Cells.Item(36, 4).Value = '''0.02309'
$ws.Cells.Item(36, 5).Value = '''  +0.84%  '
$ws.Cells.Item(37, 4).Value = '''0.06170'
$ws.Cells.Item(37, 5).Value = '''  +2.33%  '
$ws.Cells.Item(38, 4).Value = '''0.2105'
$ws.Cells.Item(38, 5).Value = '''  +0.36%  '
$ws.Cells.Item(39, 4).Value = '''0.6383'
$ws.Cells.Item(39, 5).Value = '''  +0.28%  '
$ws.Cells.Item(40, 4).Value = '''5.017'
$ws.Cells.Item(40, 5).Value = '''  +1.12%  '
$ws.Cells.Item(41, 4).Value = '''1.194'
$ws.Cells.Item(41, 5).Value = '''  -0.64%  '
$ws.Cells.Item(42, 4).Value = '''1.395'
$ws.Cells.Item(42, 5).Value = '''  +1.00%  '
$ws.Cells.Item(43, 4).Value = '''7.956'
$ws.Cells.Item(43, 5).Value = '''  +1.85%  '
$ws.Cells.Item(44, 4).Value = '''13.35'
$ws.Cells.Item(44, 5).Value = '''  +1.11%  '
$ws.Cells.Item(45, 4).Value = '''0.5949'
$ws.Cells.Item(45, 5).Value = '''  +0.56%  '
$ws.Cells.Item(46, 4).Value = '''3.744'
$ws.Cells.Item(46, 5).Value = '''  +0.93%  '
$ws.Cells.Item(47, 4).Value = '''123.29'
$ws.Cells.Item(47, 5).Value = '''  +0.75%  '
$ws.Cells.Item(48, 4).Value = '''1.976'
$ws.Cells.Item(48, 5).Value = '''  +1.27%  '
$ws.Cells.Item(49, 4).Value = '''0.06969'
$ws.Cells.Item(49, 5).Value = '''  +1.66%  '
$ws.Cells.Item(50, 4).Value = '''1.148'
$ws.Cells.Item(50, 5).Value = '''  +0.22%  '
$ws.Cells.Item(51, 4).Value = '''73.19'
$ws.Cells.Item(51, 5).Value = '''  +0.91%  '
